$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.34
$ws.Range("G2").Value = 1.42
$ws.Range("H2").Value = 9.4
$ws.Range("I2").Value = 12
$ws.Range("K2").Value = 6.2
$ws.Range("N2").Value = 5.2
$ws.Range("P2").Value = 2.42
$ws.Range("S2").Value = 2.42
$ws.Range("T2").Value = 1.87
$ws.Range("V2").Value = 1.09
$ws.Range("X2").Value = 29
$ws.Range("Y2").Value = 60
$ws.Range("Z2").Value = 980
$ws.Range("AA2").Value = 1000
$ws.Range("AC2").Value = 13.5
$ws.Range("AD2").Value = 95
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 9.800000000000001
$ws.Range("AG2").Value = 11.5
$ws.Range("AH2").Value = 44
$ws.Range("AI2").Value = 290
$ws.Range("AJ2").Value = 12
$ws.Range("AL2").Value = 55
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 5.3
$ws.Range("AO2").Value = 1000
# Row 3
$ws.Range("F3").Value = 3.6
$ws.Range("H3").Value = 2.24
$ws.Range("I3").Value = 2.38
$ws.Range("K3").Value = 3.45
$ws.Range("L3").Value = 1.5
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 3
$ws.Range("O3").Value = 1.43
$ws.Range("P3").Value = 1.68
$ws.Range("Q3").Value = 2.28
$ws.Range("V3").Value = 1.72
$ws.Range("Y3").Value = 8.4
$ws.Range("AE3").Value = 90
$ws.Range("AI3").Value = 980
# Row 4
$ws.Range("F4").Value = 2.12
$ws.Range("G4").Value = 2.14
$ws.Range("H4").Value = 3.95
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 3.6
$ws.Range("K4").Value = 3.65
$ws.Range("M4").Value = 1.06
$ws.Range("O4").Value = 1.28
$ws.Range("S4").Value = 3.25
$ws.Range("T4").Value = 1.6
$ws.Range("U4").Value = 2.04
$ws.Range("V4").Value = 1.33
$ws.Range("W4").Value = 1.88
$ws.Range("X4").Value = 18
$ws.Range("Y4").Value = 30
$ws.Range("Z4").Value = 65
$ws.Range("AA4").Value = 75
$ws.Range("AD4").Value = 15
$ws.Range("AE4").Value = 110
$ws.Range("AH4").Value = 17
$ws.Range("AI4").Value = 120
$ws.Range("AK4").Value = 21
$ws.Range("AN4").Value = 21
$ws.Range("AO4").Value = 980
# Row 5
$ws.Range("G5").Value = 15
$ws.Range("I5").Value = 1.41
$ws.Range("J5").Value = 4.9
$ws.Range("K5").Value = 5.7
$ws.Range("L5").Value = 1.34
$ws.Range("N5").Value = 4.2
$ws.Range("P5").Value = 2.12
$ws.Range("Q5").Value = 1.73
$ws.Range("R5").Value = 1.44
$ws.Range("S5").Value = 2.82
$ws.Range("V5").Value = 3.4
$ws.Range("Z5").Value = 17.5
# Row 6
$ws.Range("J6").Value = 3.1
$ws.Range("K6").Value = 3.35
$ws.Range("L6").Value = 1.51
$ws.Range("N6").Value = 3.15
$ws.Range("Q6").Value = 2.22
$ws.Range("R6").Value = 1.27
$ws.Range("T6").Value = 1.86
$ws.Range("W6").Value = 1.72
$ws.Range("AB6").Value = 980
$ws.Range("AJ6").Value = 85
$ws.Range("AK6").Value = 980
$ws.Range("AO6").Value = 600
# Row 7
$ws.Range("F7").Value = 2.58
$ws.Range("H7").Value = 2.62
$ws.Range("I7").Value = 3.05
$ws.Range("J7").Value = 3.15
$ws.Range("L7").Value = 1.5
$ws.Range("N7").Value = 2.76
$ws.Range("O7").Value = 1.44
$ws.Range("Q7").Value = 2.28
$ws.Range("R7").Value = 1.22
$ws.Range("S7").Value = 4.8
$ws.Range("W7").Value = 1.49
$ws.Range("X7").Value = 970
# Row 8
$ws.Range("F8").Value = 3.05
$ws.Range("G8").Value = 3.25
$ws.Range("J8").Value = 3.25
$ws.Range("L8").Value = 1.47
$ws.Range("N8").Value = 3.2
$ws.Range("Q8").Value = 2.16
$ws.Range("W8").Value = 1.45
# Row 9
$ws.Range("F9").Value = 2.2
$ws.Range("G9").Value = 2.42
$ws.Range("L9").Value = 1.43
$ws.Range("N9").Value = 3.5
$ws.Range("Q9").Value = 2.02
$ws.Range("W9").Value = 1.71
$ws.Range("AF9").Value = 32
# Row 10
$ws.Range("L10").Value = 1.34
$ws.Range("Q10").Value = 1.71
# Row 11
$ws.Range("F11").Value = 2.22
$ws.Range("G11").Value = 2.36
$ws.Range("L11").Value = 1.33
$ws.Range("N11").Value = 4.5
$ws.Range("P11").Value = 2.28
$ws.Range("R11").Value = 1.5
$ws.Range("S11").Value = 2.72
$ws.Range("W11").Value = 1.73
$ws.Range("AC11").Value = 14
# Row 12
$ws.Range("F12").Value = 5.3
$ws.Range("I12").Value = 1.7
$ws.Range("S12").Value = 2.84
$ws.Range("V12").Value = 2.42
$ws.Range("X12").Value = 21
# Row 13
$ws.Range("L13").Value = 1.43
$ws.Range("T13").Value = 1.78
$ws.Range("W13").Value = 1.82
$ws.Range("Y13").Value = 14.5
$ws.Range("AJ13").Value = 27
# Row 14
$ws.Range("G14").Value = 1.88
$ws.Range("L14").Value = 1.24
$ws.Range("N14").Value = 7.6
$ws.Range("AH14").Value = 14
$ws.Range("AL14").Value = 22
$ws.Range("AM14").Value = 42
# Row 15
$ws.Range("K15").Value = 4.3
$ws.Range("Q15").Value = 1.69
$ws.Range("R15").Value = 1.56
$ws.Range("S15").Value = 2.74
$ws.Range("T15").Value = 1.7
$ws.Range("W15").Value = 2.3
$ws.Range("AK15").Value = 16
$ws.Range("AO15").Value = 55
# Row 16
$ws.Range("F16").Value = 1.53
$ws.Range("G16").Value = 1.62
$ws.Range("H16").Value = 7
$ws.Range("I16").Value = 8.800000000000001
$ws.Range("J16").Value = 3.95
$ws.Range("K16").Value = 4.6
$ws.Range("L16").Value = 1.44
$ws.Range("N16").Value = 3.25
$ws.Range("P16").Value = 1.83
$ws.Range("R16").Value = 1.31
$ws.Range("U16").Value = 1.77
$ws.Range("W16").Value = 2.6
$ws.Range("Z16").Value = 1000
$ws.Range("AB16").Value = 7.2
$ws.Range("AE16").Value = 1000
$ws.Range("AI16").Value = 1000
$ws.Range("AM16").Value = 580
# Row 17
$ws.Range("H17").Value = 2.74
$ws.Range("I17").Value = 2.8
$ws.Range("J17").Value = 3.85
$ws.Range("K17").Value = 4.1
$ws.Range("L17").Value = 1.48
$ws.Range("M17").Value = 1.08
$ws.Range("N17").Value = 3.15
$ws.Range("O17").Value = 1.42
$ws.Range("P17").Value = 1.75
$ws.Range("Q17").Value = 2.22
$ws.Range("R17").Value = 1.27
$ws.Range("S17").Value = 4.3
$ws.Range("T17").Value = 1.97
$ws.Range("U17").Value = 1.96
$ws.Range("V17").Value = 1.55
$ws.Range("X17").Value = 13.5
$ws.Range("Y17").Value = 9.6
$ws.Range("Z17").Value = 17.5
$ws.Range("AA17").Value = 60
$ws.Range("AB17").Value = 9.199999999999999
$ws.Range("AC17").Value = 9.4
$ws.Range("AD17").Value = 13
$ws.Range("AE17").Value = 110
$ws.Range("AG17").Value = 12
$ws.Range("AI17").Value = 95
$ws.Range("AJ17").Value = 110
$ws.Range("AK17").Value = 34
$ws.Range("AN17").Value = 50
$ws.Range("AO17").Value = 1000
# Row 18
$ws.Range("F18").Value = 1.61
$ws.Range("G18").Value = 1.65
$ws.Range("H18").Value = 5.4
$ws.Range("I18").Value = 6.6
$ws.Range("K18").Value = 5.3
$ws.Range("L18").Value = 1.3
$ws.Range("N18").Value = 5.3
$ws.Range("O18").Value = 1.19
$ws.Range("P18").Value = 2.56
$ws.Range("Q18").Value = 1.59
$ws.Range("R18").Value = 1.6
$ws.Range("S18").Value = 2.48
$ws.Range("T18").Value = 1.62
$ws.Range("U18").Value = 2.24
$ws.Range("W18").Value = 2.52
$ws.Range("X18").Value = 28
$ws.Range("Y18").Value = 29
$ws.Range("AB18").Value = 13
$ws.Range("AD18").Value = 24
$ws.Range("AI18").Value = 260
$ws.Range("AJ18").Value = 17
$ws.Range("AK18").Value = 16
$ws.Range("AL18").Value = 70
$ws.Range("AN18").Value = 7
# Row 19
$ws.Range("H19").Value = 3.6
$ws.Range("I19").Value = 3.75
$ws.Range("Q19").Value = 1.96
$ws.Range("S19").Value = 3.4
$ws.Range("V19").Value = 1.36
$ws.Range("W19").Value = 1.81
$ws.Range("Y19").Value = 15
$ws.Range("AJ19").Value = 28
# Row 20
$ws.Range("F20").Value = 1.96
$ws.Range("G20").Value = 2.02
$ws.Range("H20").Value = 3.7
$ws.Range("N20").Value = 5.4
$ws.Range("P20").Value = 2.48
$ws.Range("Q20").Value = 1.63
$ws.Range("R20").Value = 1.59
$ws.Range("S20").Value = 2.62
$ws.Range("V20").Value = 1.32
$ws.Range("W20").Value = 1.98
# Row 21
$ws.Range("H21").Value = 2.74
$ws.Range("M21").Value = 1.1
$ws.Range("P21").Value = 1.66
$ws.Range("T21").Value = 1.98
$ws.Range("U21").Value = 2
$ws.Range("AO21").Value = 36
